$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1956388.6
$ws.Range("I76").Value = 3348945.5
$ws.Range("J76").Value = 6809
$ws.Range("K76").Value = 3348945.5
$ws.Range("L76").Value = 6809
$ws.Range("M76").Value = -3348630.5
$ws.Range("N76").Value = -7439
$ws.Range("H79").Value = 1956388.6
$ws.Range("I79").Value = 3348945.5
$ws.Range("J79").Value = 6809
$ws.Range("K79").Value = 3348945.5
$ws.Range("L79").Value = 6809
$ws.Range("M79").Value = -3347853.5
$ws.Range("N79").Value = -8993
$ws.Range("H100").Value = 920.7143
$ws.Range("I100").Value = 722.3077
$ws.Range("K100").Value = 722.3077
$ws.Range("M100").Value = -181.3077
$ws.Range("H137").Value = 1552.6923
$ws.Range("I137").Value = 1418.6
$ws.Range("K137").Value = 4255.799999999999
$ws.Range("M137").Value = -1705.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5814677.5
$ws.Range("I2").Value = 11628406
$ws.Range("J2").Value = 949
$ws.Range("K2").Value = 11628406
$ws.Range("L2").Value = 949
$ws.Range("M2").Value = -11628293
$ws.Range("N2").Value = -1175
$ws.Range("H45").Value = 1515.8125
$ws.Range("I45").Value = 1243.75
$ws.Range("J45").Value = 1787.875
$ws.Range("K45").Value = 1243.75
$ws.Range("L45").Value = 1787.875
$ws.Range("M45").Value = -866.75
$ws.Range("N45").Value = -2541.875
$ws.Range("H61").Value = 3176.7576
$ws.Range("I61").Value = 2571.4333
$ws.Range("J61").Value = 9230
$ws.Range("K61").Value = 2571.4333
$ws.Range("L61").Value = 9230
$ws.Range("M61").Value = -2359.4333
$ws.Range("N61").Value = -9654
$ws.Range("H82").Value = 63388.668
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 63388.668
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 63388.668
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -64110.668
$ws.Range("H85").Value = 63388.668
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 63388.668
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 63388.668
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -65884.66800000001
$ws.Range("H116").Value = 5814677.5
$ws.Range("I116").Value = 11628406
$ws.Range("J116").Value = 949
$ws.Range("K116").Value = 11628406
$ws.Range("L116").Value = 949
$ws.Range("M116").Value = -11626112
$ws.Range("N116").Value = -5537
$ws.Range("H122").Value = 4880
$ws.Range("I122").Value = 4880
$ws.Range("K122").Value = 14640
$ws.Range("M122").Value = -12190
$ws.Range("H136").Value = 3176.7576
$ws.Range("I136").Value = 2571.4333
$ws.Range("J136").Value = 9230
$ws.Range("K136").Value = 7714.2999
$ws.Range("L136").Value = 27690
$ws.Range("M136").Value = -5164.2999
$ws.Range("N136").Value = -32790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5814677.5
$ws.Range("I3").Value = 11628406
$ws.Range("J3").Value = 949
$ws.Range("K3").Value = 11628406
$ws.Range("L3").Value = 949
$ws.Range("M3").Value = -11628292
$ws.Range("N3").Value = -1177
$ws.Range("H20").Value = 1439.9062
$ws.Range("I20").Value = 1419.3182
$ws.Range("J20").Value = 1485.2
$ws.Range("K20").Value = 1419.3182
$ws.Range("L20").Value = 1485.2
$ws.Range("M20").Value = -1172.3182
$ws.Range("N20").Value = -1979.2
$ws.Range("H86").Value = 73096.71000000001
$ws.Range("I86").Value = 1529.9524
$ws.Range("K86").Value = 1529.9524
$ws.Range("M86").Value = -406.9523999999999
$ws.Range("H89").Value = 73096.71000000001
$ws.Range("I89").Value = 1529.9524
$ws.Range("K89").Value = 7649.762
$ws.Range("M89").Value = -2033.762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 247.81818
$ws.Range("J7").Value = 357.4
$ws.Range("L7").Value = 357.4
$ws.Range("N7").Value = -583.4
$ws.Range("H31").Value = 2371.383
$ws.Range("I31").Value = 1493.4546
$ws.Range("K31").Value = 1493.4546
$ws.Range("M31").Value = -1198.4546
$ws.Range("H34").Value = 2371.383
$ws.Range("I34").Value = 1493.4546
$ws.Range("K34").Value = 1493.4546
$ws.Range("M34").Value = -1291.4546
$ws.Range("H58").Value = 1554390.6
$ws.Range("I58").Value = 3624370.5
$ws.Range("J58").Value = 1905.6875
$ws.Range("K58").Value = 3624370.5
$ws.Range("L58").Value = 1905.6875
$ws.Range("M58").Value = -3624167.5
$ws.Range("N58").Value = -2311.6875
$ws.Range("H136").Value = 1554390.6
$ws.Range("I136").Value = 3624370.5
$ws.Range("J136").Value = 1905.6875
$ws.Range("K136").Value = 10873111.5
$ws.Range("L136").Value = 5717.0625
$ws.Range("M136").Value = -10870561.5
$ws.Range("N136").Value = -10817.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3544.9
$ws.Range("J104").Value = 3716.5557
$ws.Range("L104").Value = 11149.6671
$ws.Range("N104").Value = -16391.6671
$ws.Range("H129").Value = 49282.934
$ws.Range("J129").Value = 73455.89999999999
$ws.Range("L129").Value = 220367.7
$ws.Range("N129").Value = -230367.7
$ws.Range("H131").Value = 10101.76
$ws.Range("J131").Value = 11078.029
$ws.Range("L131").Value = 33234.087
$ws.Range("N131").Value = -43314.087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 988655.3
$ws.Range("I132").Value = 1426406.6
$ws.Range("J132").Value = 3714.9167
$ws.Range("K132").Value = 4279219.800000001
$ws.Range("L132").Value = 11144.7501
$ws.Range("M132").Value = -4276689.800000001
$ws.Range("N132").Value = -16204.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3515.3157
$ws.Range("I40").Value = 1276.3077
$ws.Range("J40").Value = 8366.5
$ws.Range("K40").Value = 1276.3077
$ws.Range("L40").Value = 1276.3077
$ws.Range("M40").Value = -1140.3077
$ws.Range("N40").Value = -8638.5
$ws.Range("H136").Value = 3007.353
$ws.Range("I136").Value = 1412.5
$ws.Range("K136").Value = 4237.5
$ws.Range("M136").Value = -1687.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 183639.72
$ws.Range("I122").Value = 183639.72
$ws.Range("K122").Value = 550919.16
$ws.Range("M122").Value = -548469.16
$ws.Range("H126").Value = 18440.3
$ws.Range("I126").Value = 24067.166
$ws.Range("K126").Value = 72201.49800000001
$ws.Range("M126").Value = -69731.49800000001
$ws.Range("H132").Value = 1877.8148
$ws.Range("I132").Value = 1465.3478
$ws.Range("K132").Value = 4396.0434
$ws.Range("M132").Value = -1866.0434

Write-Host "All edits applied"